$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(4512.4, 5012.4),
    @(5475.3, 8451.200000000001),
    @(8800, 8900),
    @(9000, 11000),
    @(9000.6, 15060.1),
    @(4800, 9220),
    @(7000, 4650.08)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r++
}
